$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# Insert 11 new rows before the current row 29 (UNION_INDICATORS) to host the
# new high-priority indicators. Excel copies the formatting of the row above
# (row 28, an "INDICATOR_NN" row) onto each newly inserted row.
$ws.Rows("29:39").Insert()

$newIndicators = @(
    "INDICATOR_24",
    "INDICATOR_25",
    "INDICATOR_26",
    "INDICATOR_27",
    "INDICATOR_28",
    "INDICATOR_29",
    "INDICATOR_30",
    "INDICATOR_31",
    "INDICATOR_32",
    "INDICATOR_34",
    "INDICATOR_35"
)

$row = 29
foreach ($name in $newIndicators) {
    $ws.Cells.Item($row, 1).Value = "CREATE/MODIFY"
    $ws.Cells.Item($row, 2).Value = "LIB_EWS_IT"
    $ws.Cells.Item($row, 3).Value = $name
    $ws.Cells.Item($row, 5).Value = "String"
    $ws.Cells.Item($row, 6).Value = "String"
    $row = $row + 1
}

# Rows 38 and 39 (INDICATOR_34 / INDICATOR_35) keep the plain "s=1" look
# instead of the "s=5" formatting inherited from the row above - match that
# by copying the format already used by the similar C3 cell.
$ws.Range("C3").Copy()
$ws.Range("C38").PasteSpecial(-4122)
$ws.Range("C39").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the view: after editing, Excel had scrolled back up and selected E36
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 27
$ws.Range("E36").Select() | Out-Null
